$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: rename rec_size -> my_size (I1), model_group -> rec_size (J1)
$ws.Range("I1").Value = "my_size"
$ws.Range("J1").Value = "rec_size"

# Row 2 updates
$ws.Range("E2").Value = "KSP Almaty-1"
$ws.Range("I2").Value = "XL"
$ws.Range("J2").Value = "XL"

# Row 3 updates
$ws.Range("E3").Value = "KSP Astana-2"
$ws.Range("J3").Value = "M"
